$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 22222624
$ws.Range("I9").Value = 27777876
$ws.Range("K9").Value = 27777876
$ws.Range("M9").Value = -27777707
$ws.Range("H15").Value = 639.5294
$ws.Range("I15").Value = 639.5294
$ws.Range("K15").Value = 1918.5882
$ws.Range("M15").Value = -1749.5882
$ws.Range("H33").Value = 391.44446
$ws.Range("I33").Value = 253.28572
$ws.Range("J33").Value = 875
$ws.Range("K33").Value = 253.28572
$ws.Range("L33").Value = 875
$ws.Range("M33").Value = -24.28572
$ws.Range("N33").Value = -1333
$ws.Range("H53").Value = 1733.3334
$ws.Range("J53").Value = 1733.3334
$ws.Range("L53").Value = 1733.3334
$ws.Range("N53").Value = -3007.3334
$ws.Range("H88").Value = 632329.8
$ws.Range("J88").Value = 68649.53
$ws.Range("L88").Value = 68649.53
$ws.Range("N88").Value = -69461.53
$ws.Range("H91").Value = 632329.8
$ws.Range("J91").Value = 68649.53
$ws.Range("L91").Value = 68649.53
$ws.Range("N91").Value = -71457.53
$ws.Range("H132").Value = 1558.6111
$ws.Range("J132").Value = 1412
$ws.Range("L132").Value = 4236
$ws.Range("N132").Value = -9296
$ws.Range("H137").Value = 428878
$ws.Range("I137").Value = 1743.6364
$ws.Range("K137").Value = 5230.9092
$ws.Range("M137").Value = -2680.9092
$ws.Range("H138").Value = 1525.4546
$ws.Range("J138").Value = 2498.8572
$ws.Range("L138").Value = 7496.571599999999
$ws.Range("N138").Value = -17776.5716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 666666.7
$ws.Range("I13").Value = 666666.7
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 666666.7
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -666522.7
$ws.Range("N13").ClearContents()
$ws.Range("H15").Value = 4998.75
$ws.Range("J15").Value = 4998.75
$ws.Range("L15").Value = 4998.75
$ws.Range("N15").Value = -5698.75
$ws.Range("H35").Value = 1863
$ws.Range("I35").Value = 1863
$ws.Range("K35").Value = 1863
$ws.Range("M35").Value = -1457
$ws.Range("H61").Value = 128001.125
$ws.Range("I61").Value = 3585
$ws.Range("J61").Value = 501249.5
$ws.Range("K61").Value = 3585
$ws.Range("L61").Value = 501249.5
$ws.Range("M61").Value = -3373
$ws.Range("N61").Value = -501673.5
$ws.Range("H132").Value = 2466.8215
$ws.Range("I132").Value = 2263.9565
$ws.Range("K132").Value = 6791.869499999999
$ws.Range("M132").Value = -4261.869499999999
$ws.Range("H136").Value = 128001.125
$ws.Range("I136").Value = 3585
$ws.Range("J136").Value = 501249.5
$ws.Range("K136").Value = 10755
$ws.Range("L136").Value = 1503748.5
$ws.Range("M136").Value = -8205
$ws.Range("N136").Value = -1508848.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 103140.09
$ws.Range("I20").Value = 136532.42
$ws.Range("J20").Value = 2963.125
$ws.Range("K20").Value = 136532.42
$ws.Range("L20").Value = 2963.125
$ws.Range("M20").Value = -136285.42
$ws.Range("N20").Value = -3457.125
$ws.Range("H25").Value = 2875.1667
$ws.Range("I25").Value = 2050.2
$ws.Range("J25").Value = 7000
$ws.Range("K25").Value = 2050.2
$ws.Range("L25").Value = 7000
$ws.Range("M25").Value = -1815.2
$ws.Range("N25").Value = -7470
$ws.Range("H105").Value = 64973.625
$ws.Range("I105").Value = 101392.9
$ws.Range("J105").Value = 4274.8335
$ws.Range("K105").Value = 101392.9
$ws.Range("L105").Value = 4274.8335
$ws.Range("M105").Value = -99645.89999999999
$ws.Range("N105").Value = -7768.8335
$ws.Range("H132").Value = 29420.525
$ws.Range("J132").Value = 29420.525
$ws.Range("L132").Value = 29420.525
$ws.Range("N132").Value = -39540.525

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 718
$ws.Range("J14").Value = 900
$ws.Range("L14").Value = 900
$ws.Range("N14").Value = -1240
$ws.Range("H31").Value = 3192.0938
$ws.Range("I31").Value = 1995.9445
$ws.Range("J31").Value = 4730
$ws.Range("K31").Value = 1995.9445
$ws.Range("L31").Value = 4730
$ws.Range("M31").Value = -1700.9445
$ws.Range("N31").Value = -5320
$ws.Range("H34").Value = 3192.0938
$ws.Range("I34").Value = 1995.9445
$ws.Range("J34").Value = 4730
$ws.Range("K34").Value = 1995.9445
$ws.Range("L34").Value = 4730
$ws.Range("M34").Value = -1793.9445
$ws.Range("N34").Value = -5134
$ws.Range("H132").Value = 543387.7
$ws.Range("I132").Value = 1820.091
$ws.Range("K132").Value = 5460.272999999999
$ws.Range("M132").Value = -2930.272999999999
$ws.Range("H134").Value = 94698.45
$ws.Range("I134").Value = 4405.5
$ws.Range("K134").Value = 13216.5
$ws.Range("M134").Value = -10681.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 326.8095
$ws.Range("I38").Value = 24.166666
$ws.Range("J38").Value = 447.86667
$ws.Range("K38").Value = 72.49999800000001
$ws.Range("L38").Value = 1343.60001
$ws.Range("M38").Value = 274.500002
$ws.Range("N38").Value = -2037.60001
$ws.Range("H140").Value = 4175
$ws.Range("I140").Value = 4175
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 12525
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -7345
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 9599.875
$ws.Range("I43").Value = 5971.2856
$ws.Range("J43").Value = 35000
$ws.Range("K43").Value = 5971.2856
$ws.Range("L43").Value = 35000
$ws.Range("M43").Value = -5820.2856
$ws.Range("N43").Value = -35302
$ws.Range("H48").Value = 6000
$ws.Range("I48").Value = 6000
$ws.Range("K48").Value = 6000
$ws.Range("M48").Value = -5515
$ws.Range("H52").Value = 18815.5
$ws.Range("J52").Value = 18578.6
$ws.Range("L52").Value = 18578.6
$ws.Range("N52").Value = -19096.6
$ws.Range("H122").Value = 7016121.5
$ws.Range("I122").Value = 8018313.5
$ws.Range("K122").Value = 24054940.5
$ws.Range("M122").Value = -24052490.5
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H132").Value = 4191.4585
$ws.Range("I132").Value = 3684.9285
$ws.Range("J132").Value = 4900.6
$ws.Range("K132").Value = 11054.7855
$ws.Range("L132").Value = 14701.8
$ws.Range("M132").Value = -8524.7855
$ws.Range("N132").Value = -19761.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3248.1
$ws.Range("I7").Value = 2415.1667
$ws.Range("J7").Value = 4497.5
$ws.Range("K7").Value = 2415.1667
$ws.Range("L7").Value = 4497.5
$ws.Range("M7").Value = -2303.1667
$ws.Range("N7").Value = -4721.5
$ws.Range("H24").Value = 25000
$ws.Range("J24").Value = 25000
$ws.Range("L24").Value = 25000
$ws.Range("N24").Value = -25686
$ws.Range("H35").Value = 1096.7778
$ws.Range("I35").Value = 1096.7778
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1096.7778
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -760.7778000000001
$ws.Range("N35").ClearContents()
$ws.Range("H46").Value = 1587
$ws.Range("I46").Value = 1584.4
$ws.Range("J46").Value = 1600
$ws.Range("K46").Value = 1584.4
$ws.Range("L46").Value = 1600
$ws.Range("M46").Value = -1396.4
$ws.Range("N46").Value = -1976
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H100").Value = 4577.5
$ws.Range("I100").Value = 3992.2222
$ws.Range("K100").Value = 3992.2222
$ws.Range("M100").Value = -3451.2222
$ws.Range("H126").Value = 3248.1
$ws.Range("I126").Value = 2415.1667
$ws.Range("J126").Value = 4497.5
$ws.Range("K126").Value = 7245.500100000001
$ws.Range("L126").Value = 13492.5
$ws.Range("M126").Value = -4775.500100000001
$ws.Range("N126").Value = -18432.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H42").Value = 15000
$ws.Range("J42").Value = 15000
$ws.Range("L42").Value = 15000
$ws.Range("N42").Value = -15756
$ws.Range("H43").Value = 6001
$ws.Range("I43").Value = 6001
$ws.Range("K43").Value = 6001
$ws.Range("M43").Value = -5852
$ws.Range("H107").Value = 8287.9375
$ws.Range("I107").Value = 12778.223
$ws.Range("K107").Value = 38334.669
$ws.Range("M107").Value = -36414.669
$ws.Range("H132").Value = 1828.0217
$ws.Range("I132").Value = 1532.2812
$ws.Range("K132").Value = 4596.8436
$ws.Range("M132").Value = -2066.8436
$ws.Range("H136").Value = 1465
$ws.Range("I136").Value = 1363.25
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 4089.75
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -1539.75
$ws.Range("N136").Value = -15600
